$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value2 = 44656
$ws.Range("J3").Value2 = 85
$ws.Range("K3").Value2 = 5000
$ws.Range("L3").Value2 = 5000
$ws.Range("M3").Value2 = 5000
$ws.Range("P3").Value2 = 5000

# Row 5
$ws.Range("D5").Value2 = 44649
$ws.Range("J5").Value2 = 20

# Row 6
$ws.Range("D6").Value2 = 44315
$ws.Range("J6").Value2 = 40
$ws.Range("K6").Value2 = 4000
$ws.Range("L6").Value2 = 4000
$ws.Range("M6").Value2 = 4000
$ws.Range("P6").Value2 = 4000

# Row 7
$ws.Range("D7").Value2 = 44749
$ws.Range("J7").Value2 = 65
$ws.Range("K7").Value2 = 6000
$ws.Range("L7").Value2 = 6000
$ws.Range("M7").Value2 = 6000
$ws.Range("P7").Value2 = 6000

# Row 8
$ws.Range("D8").Value2 = 44956
$ws.Range("J8").Value2 = 40

# Row 9
$ws.Range("D9").Value2 = 44966
$ws.Range("J9").Value2 = 40

# Row 10
$ws.Range("D10").Value2 = 44365
$ws.Range("J10").Value2 = 55
$ws.Range("K10").Value2 = 5000
$ws.Range("L10").Value2 = 5000
$ws.Range("M10").Value2 = 5000
$ws.Range("P10").Value2 = 5000

# Row 11
$ws.Range("D11").Value2 = 44957
$ws.Range("J11").Value2 = 20

# Row 12
$ws.Range("D12").Value2 = 44498
$ws.Range("J12").Value2 = 40
$ws.Range("K12").Value2 = 4000
$ws.Range("L12").Value2 = 4000
$ws.Range("M12").Value2 = 4000
$ws.Range("P12").Value2 = 4000

# Row 13
$ws.Range("D13").Value2 = 44959
$ws.Range("J13").Value2 = 40
$ws.Range("K13").Value2 = 5000
$ws.Range("L13").Value2 = 5000
$ws.Range("M13").Value2 = 5000
$ws.Range("P13").Value2 = 5000

# Row 14
$ws.Range("D14").Value2 = 44508
$ws.Range("J14").Value2 = 30

# Row 15
$ws.Range("D15").Value2 = 44291
$ws.Range("J15").Value2 = 35
$ws.Range("K15").Value2 = 4000
$ws.Range("L15").Value2 = 4000
$ws.Range("M15").Value2 = 4000
$ws.Range("P15").Value2 = 4000

# Row 16
$ws.Range("D16").Value2 = 44497
$ws.Range("J16").Value2 = 20

# Row 17
$ws.Range("D17").Value2 = 44390
$ws.Range("J17").Value2 = 55
$ws.Range("K17").Value2 = 6000
$ws.Range("L17").Value2 = 6000
$ws.Range("M17").Value2 = 6000
$ws.Range("P17").Value2 = 6000

# Row 18
$ws.Range("D18").Value2 = 44313
$ws.Range("J18").Value2 = 20
$ws.Range("K18").Value2 = 4000
$ws.Range("L18").Value2 = 4000
$ws.Range("M18").Value2 = 4000
$ws.Range("P18").Value2 = 4000

# Row 19
$ws.Range("D19").Value2 = 44680
$ws.Range("K19").Value2 = 5000
$ws.Range("L19").Value2 = 5000
$ws.Range("M19").Value2 = 5000
$ws.Range("P19").Value2 = 5000

# Row 20
$ws.Range("D20").Value2 = 44781

# Row 21
$ws.Range("D21").Value2 = 44301
$ws.Range("J21").Value2 = 40
$ws.Range("K21").Value2 = 3000
$ws.Range("L21").Value2 = 3000
$ws.Range("M21").Value2 = 3000
$ws.Range("P21").Value2 = 3000

# Row 22
$ws.Range("D22").Value2 = 44312
$ws.Range("J22").Value2 = 50
$ws.Range("K22").Value2 = 4000
$ws.Range("L22").Value2 = 4000
$ws.Range("M22").Value2 = 4000
$ws.Range("P22").Value2 = 4000

# Row 23
$ws.Range("D23").Value2 = 44679
$ws.Range("J23").Value2 = 50
$ws.Range("K23").Value2 = 5000
$ws.Range("L23").Value2 = 5000
$ws.Range("M23").Value2 = 5000
$ws.Range("P23").Value2 = 5000

# Row 24
$ws.Range("D24").Value2 = 44259

# Row 25
$ws.Range("D25").Value2 = 44504
$ws.Range("J25").Value2 = 55

# Row 26
$ws.Range("D26").Value2 = 44777
$ws.Range("J26").Value2 = 25
$ws.Range("K26").Value2 = 5000
$ws.Range("L26").Value2 = 5000
$ws.Range("M26").Value2 = 5000
$ws.Range("P26").Value2 = 5000

# Row 27
$ws.Range("D27").Value2 = 44316
$ws.Range("K27").Value2 = 4000
$ws.Range("L27").Value2 = 4000
$ws.Range("M27").Value2 = 4000
$ws.Range("P27").Value2 = 4000

# Row 28
$ws.Range("D28").Value2 = 44176
$ws.Range("J28").Value2 = 10
